$wb = $excel.ActiveWorkbook

# ---- Sheet "hpi" ----
$ws = $wb.Worksheets.Item("hpi")
$ws.Range("C2").Value = "Antacids no longer providing relief is present."
$ws.Range("D2").Value = "This suggests a non-Scleroderma related cause for the symptoms, as Scleroderma typically involves esophageal dysmotility that may respond to antacids."

$ws.Range("C3").Value = "Long-standing heartburn (duration of years) is present."
$ws.Range("D3").Value = "Chronic heartburn is more commonly associated with gastroesophageal reflux disease (GERD) rather than Scleroderma."

$ws.Range("C4").Value = "Long-standing reflux (duration of years) is present."
$ws.Range("D4").Value = "Chronic reflux symptoms are more indicative of GERD and less specific to Scleroderma."

$ws.Range("C5").Value = "Pain relieved with regurgitation is present."
$ws.Range("D5").Value = "This finding suggests a mechanical issue related to reflux rather than esophageal motility problems seen in Scleroderma."

$ws.Range("C6").Value = "Pain worse when lying down (positional) is present."
$ws.Range("D6").Value = "Positional pain is more characteristic of reflux conditions rather than the esophageal involvement seen in Scleroderma."

# ---- Sheet "hist" ----
$ws = $wb.Worksheets.Item("hist")
$ws.Range("C2").Value = "Absence of prior treatment with radiation to the neck, arm, or jaw"
$ws.Range("D2").Value = "Radiation exposure is a known risk factor for scleroderma; its absence weakens the likelihood of the diagnosis."

$ws.Range("C3").Value = "Absence of diagnosed hypertension"
$ws.Range("D3").Value = "Hypertension is often associated with scleroderma; its absence suggests a lower likelihood of the disease."

$ws.Range("A4").Value = "Amlodipine"
$ws.Range("B4").Value = "Amlodipine is a calcium channel blocker that may be used to manage hypertension, which can be associated with scleroderma-related complications."
$ws.Range("C4").Value = "Absence of diagnosed coronary artery disease"
$ws.Range("D4").Value = "Coronary artery disease can be a complication of scleroderma; its absence may indicate a lower risk for the condition."

$ws.Range("A5").Value = "Absence of alcohol use disorder"
$ws.Range("B5").Value = "While not directly indicative, the absence of alcohol use disorder may suggest a lower likelihood of other conditions that could mimic scleroderma."
$ws.Range("C5").Value = "Absence of prior myocardial infarction"
$ws.Range("D5").Value = "A history of myocardial infarction is often seen in patients with scleroderma; its absence suggests a lower likelihood of the diagnosis."

$ws.Range("A6").Value = "Absence of nicotine dependence"
$ws.Range("B6").Value = "Similar to alcohol use disorder, the absence of nicotine dependence may indicate a lower risk for conditions that could complicate or mimic scleroderma."
$ws.Range("C6").Value = "Absence of obesity"
$ws.Range("D6").Value = "Obesity can be a risk factor for various conditions, including those that may mimic scleroderma; its absence may indicate a lower risk."

# ---- Sheet "soc" ----
$ws = $wb.Worksheets.Item("soc")
$ws.Range("C2").Value = "Family history of cancer is absent."
$ws.Range("D2").Value = "A lack of family history of cancer may suggest a lower risk for certain autoimmune diseases, including Scleroderma."

$ws.Range("C3").Value = "Recent Travel is absent."
$ws.Range("D3").Value = "Absence of recent travel may indicate lower exposure to environmental triggers that can exacerbate autoimmune conditions."

$ws.Range("C4").Value = "Recent medical procedure is absent."
$ws.Range("D4").Value = "Absence of recent medical procedures may suggest a lower likelihood of complications or triggers related to Scleroderma."

$ws.Range("A5").Value = "Family history of Rheumatoid Arthritis is absent."
$ws.Range("B5").Value = "While this finding does not directly support Scleroderma, the absence of a common autoimmune condition may suggest a different autoimmune profile."
$ws.Range("C5").Value = "Gestational complications with prior pregnancy is absent."
$ws.Range("D5").Value = "Absence of gestational complications may indicate a lower risk for autoimmune conditions that can be influenced by pregnancy."

$ws.Range("A6").Value = "Alcohol use is absent."
$ws.Range("B6").Value = "Absence of alcohol use may indicate a healthier lifestyle, which can be a factor in managing autoimmune diseases."
$ws.Range("C6").Value = "Current tobacco use is absent."
$ws.Range("D6").Value = "Current non-use of tobacco may suggest a lower risk for developing autoimmune diseases, including Scleroderma."

# ---- Sheet "obj" ----
$ws = $wb.Worksheets.Item("obj")
$ws.Range("D2").Value = "Hand thickening is a classic sign of Scleroderma; its absence strongly suggests that the diagnosis is unlikely."
$ws.Range("D3").Value = "Finger ulcers are a common manifestation of Scleroderma; their absence is a strong indicator against the diagnosis."
$ws.Range("D4").Value = "Muscle weakness can occur in Scleroderma; its absence suggests that the diagnosis is less likely."
$ws.Range("B5").Value = "A hoarse voice can be associated with esophageal involvement in Scleroderma; its absence does not support the diagnosis but is not definitive."
$ws.Range("D5").Value = "Joint swelling can be associated with Scleroderma; its absence is a strong indicator against the diagnosis."
$ws.Range("B6").Value = "Cough can be associated with pulmonary involvement in Scleroderma; its absence does not support the diagnosis but is not definitive."
$ws.Range("D6").Value = "Rheumatoid nodules are not typically associated with Scleroderma; their absence does not support the diagnosis."
